$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s-val data (regen sval data to filter save games)
$data = @(
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 0, 28.30127388105354),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 0, 11.47044854674929),
    @(1.505614041169197, 1.65323645889881, 157.8057217802531, 6.48142807727062, 0, 167.4460003575917),
    @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.594575437922795),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
